# Regenerate merged AHB files
#
# 1) Rename the header-row labels from the "_old"/"_new" naming scheme to the
#    versioned "_FV2404"/"_FV2410" naming scheme (the "diff" column header is
#    unchanged).
# 2) Freeze the header row (top row) so it stays visible while scrolling.
# 3) Turn the A1:U92 data range into a native Excel Table ("Table1") with an
#    autofilter, so the header labels above become the table's column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header-row labels (row 1, columns A..U) ---------------------
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2) Freeze the top (header) row -----------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true | Out-Null

# --- 3) Convert the used range into a native Table --------------------------
$dataRange = $ws.Range("A1:U92")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
